$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Split the generic "F1" / "F2" cross labels into distinct sub-labels
$ws.Range("A4").Value = "F1a"
$ws.Range("A5").Value = "F1b"
$ws.Range("A6").Value = "F2a"
$ws.Range("A7").Value = "F2b"

# Leave the selection on E7, matching the saved workbook state
$ws.Range("E7").Select()
